$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl4"
$ws.Range("C2").Value = "Ccr3"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.8774723333333333
$ws.Range("H2").Value = 2.632417
$ws.Range("I2").Value = 0.0007088393434259271
$ws.Range("J2").Value = 0.0007088393434259271
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1808983333333333
$ws.Range("N2").Value = 0.5426949999999999
$ws.Range("O2").Value = 0.09546831801815302
$ws.Range("P2").Value = 0.09546831801815302
$ws.Range("Q2").Value = 0.1587332826461111
$ws.Range("R2").Value = 1.428599543815
$ws.Range("S2").Value = 0.00006767169986196518
$ws.Range("T2").Value = 0.00006767169986196518

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl4"
$ws.Range("C3").Value = "Ccr3"
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.8774723333333333
$ws.Range("H3").Value = 2.632417
$ws.Range("I3").Value = 0.0007088393434259271
$ws.Range("J3").Value = 0.0007088393434259271
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.572737
$ws.Range("N3").Value = 4.718211
$ws.Range("O3").Value = 0.8300051930177132
$ws.Range("P3").Value = 0.8300051930177132
$ws.Range("Q3").Value = 1.380033205109667
$ws.Range("R3").Value = 12.420298845987
$ws.Range("S3").Value = 0.0005883403360587857
$ws.Range("T3").Value = 0.0005883403360587857

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ccl4"
$ws.Range("C4").Value = "Ccr3"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.8774723333333333
$ws.Range("H4").Value = 2.632417
$ws.Range("I4").Value = 0.0007088393434259271
$ws.Range("J4").Value = 0.0007088393434259271
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1412166666666667
$ws.Range("N4").Value = 0.42365
$ws.Range("O4").Value = 0.0745264889641337
$ws.Range("P4").Value = 0.07452648896413369
$ws.Range("Q4").Value = 0.1239137180055556
$ws.Range("R4").Value = 1.11522346205
$ws.Range("S4").Value = 0.00005282730750517614
$ws.Range("T4").Value = 0.00005282730750517613

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl4"
$ws.Range("C5").Value = "Ccr3"
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.084971
$ws.Range("H5").Value = 0.254913
$ws.Range("I5").Value = 0.00006864123866041489
$ws.Range("J5").Value = 0.00006864123866041489
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1808983333333333
$ws.Range("N5").Value = 0.5426949999999999
$ws.Range("O5").Value = 0.09546831801815302
$ws.Range("P5").Value = 0.09546831801815302
$ws.Range("Q5").Value = 0.01537111228166666
$ws.Range("R5").Value = 0.138340010535
$ws.Range("S5").Value = 0.000006553063601592428
$ws.Range("T5").Value = 0.000006553063601592428

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ccl4"
$ws.Range("C6").Value = "Ccr3"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.084971
$ws.Range("H6").Value = 0.254913
$ws.Range("I6").Value = 0.00006864123866041489
$ws.Range("J6").Value = 0.00006864123866041489
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.572737
$ws.Range("N6").Value = 4.718211
$ws.Range("O6").Value = 0.8300051930177132
$ws.Range("P6").Value = 0.8300051930177132
$ws.Range("Q6").Value = 0.133637035627
$ws.Range("R6").Value = 1.202733320643
$ws.Range("S6").Value = 0.00005697258454331258
$ws.Range("T6").Value = 0.00005697258454331258

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ccl4"
$ws.Range("C7").Value = "Ccr3"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.084971
$ws.Range("H7").Value = 0.254913
$ws.Range("I7").Value = 0.00006864123866041489
$ws.Range("J7").Value = 0.00006864123866041489
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1412166666666667
$ws.Range("N7").Value = 0.42365
$ws.Range("O7").Value = 0.0745264889641337
$ws.Range("P7").Value = 0.07452648896413369
$ws.Range("Q7").Value = 0.01199932138333334
$ws.Range("R7").Value = 0.10799389245
$ws.Range("S7").Value = 0.000005115590515509878
$ws.Range("T7").Value = 0.000005115590515509877

$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Ccl4"
$ws.Range("C8").Value = "Ccr3"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 71.07177233333334
$ws.Range("H8").Value = 213.215317
$ws.Range("I8").Value = 0.05741317021985155
$ws.Range("J8").Value = 0.05741317021985154
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1808983333333333
$ws.Range("N8").Value = 0.5426949999999999
$ws.Range("O8").Value = 0.09546831801815302
$ws.Range("P8").Value = 0.09546831801815302
$ws.Range("Q8").Value = 12.85676516214611
$ws.Range("R8").Value = 115.710886459315
$ws.Range("S8").Value = 0.00548113879297914
$ws.Range("T8").Value = 0.005481138792979139

$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Ccl4"
$ws.Range("C9").Value = "Ccr3"
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 71.07177233333334
$ws.Range("H9").Value = 213.215317
$ws.Range("I9").Value = 0.05741317021985155
$ws.Range("J9").Value = 0.05741317021985154
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.572737
$ws.Range("N9").Value = 4.718211
$ws.Range("O9").Value = 0.8300051930177132
$ws.Range("P9").Value = 0.8300051930177132
$ws.Range("Q9").Value = 111.7772060042097
$ws.Range("R9").Value = 1005.994854037887
$ws.Range("S9").Value = 0.04765322943008671
$ws.Range("T9").Value = 0.0476532294300867

$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Ccl4"
$ws.Range("C10").Value = "Ccr3"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 71.07177233333334
$ws.Range("H10").Value = 213.215317
$ws.Range("I10").Value = 0.05741317021985155
$ws.Range("J10").Value = 0.05741317021985154
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1412166666666667
$ws.Range("N10").Value = 0.42365
$ws.Range("O10").Value = 0.0745264889641337
$ws.Range("P10").Value = 0.07452648896413369
$ws.Range("Q10").Value = 10.03651878300556
$ws.Range("R10").Value = 90.32866904705001
$ws.Range("S10").Value = 0.004278801996785696
$ws.Range("T10").Value = 0.004278801996785695

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Ccl4"
$ws.Range("C11").Value = "Ccr3"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.03596933333333333
$ws.Range("H11").Value = 0.107908
$ws.Range("I11").Value = 0.00002905673222380989
$ws.Range("J11").Value = 0.00002905673222380989
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1808983333333333
$ws.Range("N11").Value = 0.5426949999999999
$ws.Range("O11").Value = 0.09546831801815302
$ws.Range("P11").Value = 0.09546831801815302
$ws.Range("Q11").Value = 0.006506792451111109
$ws.Range("R11").Value = 0.05856113205999999
$ws.Range("S11").Value = 0.000002773997352510997
$ws.Range("T11").Value = 0.000002773997352510997

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Ccl4"
$ws.Range("C12").Value = "Ccr3"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.03596933333333333
$ws.Range("H12").Value = 0.107908
$ws.Range("I12").Value = 0.00002905673222380989
$ws.Range("J12").Value = 0.00002905673222380989
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.572737
$ws.Range("N12").Value = 4.718211
$ws.Range("O12").Value = 0.8300051930177132
$ws.Range("P12").Value = 0.8300051930177132
$ws.Range("Q12").Value = 0.05657030139866667
$ws.Range("R12").Value = 0.5091327125880001
$ws.Range("S12").Value = 0.00002411723863788733
$ws.Range("T12").Value = 0.00002411723863788733

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Ccl4"
$ws.Range("C13").Value = "Ccr3"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.03596933333333333
$ws.Range("H13").Value = 0.107908
$ws.Range("I13").Value = 0.00002905673222380989
$ws.Range("J13").Value = 0.00002905673222380989
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1412166666666667
$ws.Range("N13").Value = 0.42365
$ws.Range("O13").Value = 0.0745264889641337
$ws.Range("P13").Value = 0.07452648896413369
$ws.Range("Q13").Value = 0.005079469355555556
$ws.Range("R13").Value = 0.0457152242
$ws.Range("S13").Value = 0.000002165496233411556
$ws.Range("T13").Value = 0.000002165496233411555

$ws.Range("A14").Value = "Neutrophils"
$ws.Range("B14").Value = "Ccl4"
$ws.Range("C14").Value = "Ccr3"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1112.75885
$ws.Range("H14").Value = 3338.27655
$ws.Range("I14").Value = 0.8989084016233635
$ws.Range("J14").Value = 0.8989084016233634
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1808983333333333
$ws.Range("N14").Value = 0.5426949999999999
$ws.Range("O14").Value = 0.09546831801815302
$ws.Range("P14").Value = 0.09546831801815302
$ws.Range("Q14").Value = 201.2962213669167
$ws.Range("R14").Value = 1811.66599230225
$ws.Range("S14").Value = 0.08581727315536888
$ws.Range("T14").Value = 0.08581727315536887

$ws.Range("A15").Value = "Neutrophils"
$ws.Range("B15").Value = "Ccl4"
$ws.Range("C15").Value = "Ccr3"
$ws.Range("D15").Value = "Neutrophils"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1112.75885
$ws.Range("H15").Value = 3338.27655
$ws.Range("I15").Value = 0.8989084016233635
$ws.Range("J15").Value = 0.8989084016233634
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.572737
$ws.Range("N15").Value = 4.718211
$ws.Range("O15").Value = 0.8300051930177132
$ws.Range("P15").Value = 0.8300051930177132
$ws.Range("Q15").Value = 1750.07701547245
$ws.Range("R15").Value = 15750.69313925205
$ws.Range("S15").Value = 0.7460986413946439
$ws.Range("T15").Value = 0.7460986413946439

$ws.Range("A16").Value = "Neutrophils"
$ws.Range("B16").Value = "Ccl4"
$ws.Range("C16").Value = "Ccr3"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1112.75885
$ws.Range("H16").Value = 3338.27655
$ws.Range("I16").Value = 0.8989084016233635
$ws.Range("J16").Value = 0.8989084016233634
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.1412166666666667
$ws.Range("N16").Value = 0.42365
$ws.Range("O16").Value = 0.0745264889641337
$ws.Range("P16").Value = 0.07452648896413369
$ws.Range("Q16").Value = 157.1400956008334
$ws.Range("R16").Value = 1414.2608604075
$ws.Range("S16").Value = 0.06699248707335066
$ws.Range("T16").Value = 0.06699248707335065

$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Ccl4"
$ws.Range("C17").Value = "Ccr3"
$ws.Range("D17").Value = "Inflammatory-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 53.07112033333333
$ws.Range("H17").Value = 159.213361
$ws.Range("I17").Value = 0.04287189084247485
$ws.Range("J17").Value = 0.04287189084247484
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.1808983333333333
$ws.Range("N17").Value = 0.5426949999999999
$ws.Range("O17").Value = 0.09546831801815302
$ws.Range("P17").Value = 0.09546831801815302
$ws.Range("Q17").Value = 9.600477216432775
$ws.Range("R17").Value = 86.40429494789498
$ws.Range("S17").Value = 0.004092907308988931
$ws.Range("T17").Value = 0.004092907308988931

$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("B18").Value = "Ccl4"
$ws.Range("C18").Value = "Ccr3"
$ws.Range("D18").Value = "Neutrophils"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 53.07112033333333
$ws.Range("H18").Value = 159.213361
$ws.Range("I18").Value = 0.04287189084247485
$ws.Range("J18").Value = 0.04287189084247484
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 1.572737
$ws.Range("N18").Value = 4.718211
$ws.Range("O18").Value = 0.8300051930177132
$ws.Range("P18").Value = 0.8300051930177132
$ws.Range("Q18").Value = 83.46691457968566
$ws.Range("R18").Value = 751.202231217171
$ws.Range("S18").Value = 0.03558389203374267
$ws.Range("T18").Value = 0.03558389203374267

$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("B19").Value = "Ccl4"
$ws.Range("C19").Value = "Ccr3"
$ws.Range("D19").Value = "Resolving-Mac"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 53.07112033333333
$ws.Range("H19").Value = 159.213361
$ws.Range("I19").Value = 0.04287189084247485
$ws.Range("J19").Value = 0.04287189084247484
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0.3333333333333333
$ws.Range("M19").Value = 0.1412166666666667
$ws.Range("N19").Value = 0.42365
$ws.Range("O19").Value = 0.0745264889641337
$ws.Range("P19").Value = 0.07452648896413369
$ws.Range("Q19").Value = 7.494526709738889
$ws.Range("R19").Value = 67.45074038765
$ws.Range("S19").Value = 0.003195091499743247
$ws.Range("T19").Value = 0.003195091499743246

